# Update "想去人数" (interested-count) values on the 展览 and 全部类型 sheets
# to reflect the refreshed scrape data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 sheet updates (row => new F value)
$wsExhibit.Range("F9").Value  = 1473
$wsExhibit.Range("F11").Value = 3072
$wsExhibit.Range("F31").Value = 30
$wsExhibit.Range("F35").Value = 1637
$wsExhibit.Range("F37").Value = 1008
$wsExhibit.Range("F38").Value = 34
$wsExhibit.Range("F46").Value = 206

# 全部类型 sheet updates (same events, different row offsets)
$wsAll.Range("F10").Value = 1473
$wsAll.Range("F13").Value = 3072
$wsAll.Range("F34").Value = 30
$wsAll.Range("F36").Value = 1637
$wsAll.Range("F38").Value = 1008
$wsAll.Range("F39").Value = 34
$wsAll.Range("F47").Value = 206
